# Weekly CompStat (42nd Precinct) refresh: shift the reporting week forward
# by one week and update all the Week-to-Date / 28-Day / Year-to-Date /
# historical % columns with the newly collected crime figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: volume/week number and the "Report Covering the Week ... Through"
# banner (these live inside rich-text shared strings, but Range.Value just
# needs the plain new text).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# ---------------------------------------------------------------------------
# Helper: a couple of cells (C15 / C26) flip from a plain numeric 0/1 count
# to the sheet's "no data" text placeholder -- a literal text "0" styled
# like the other placeholder cells (e.g. C14/D14). We reproduce that by
# forcing the value to text (leading apostrophe) and then lifting the
# number format / font / alignment from an existing placeholder cell
# (C14) so the cell matches the rest of the "0" placeholders exactly.
# ---------------------------------------------------------------------------
function Set-PlaceholderZero($addr) {
    $ws.Range($addr).Value = "'0"
    $ws.Range("C14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 33.333333333333
$ws.Range("L14").Value = -55.555555555555
$ws.Range("M14").Value = -20
$ws.Range("N14").Value = -55.555555555555

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-PlaceholderZero "C15"
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 50
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = 6.25
$ws.Range("N15").Value = -32

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = -57.142857142857
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 144
$ws.Range("J16").Value = 125
$ws.Range("K16").Value = 15.2
$ws.Range("L16").Value = 29.729729729729
$ws.Range("M16").Value = 45.454545454545
$ws.Range("N16").Value = -67.713004484304

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 12
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 45
$ws.Range("G17").Value = 52
$ws.Range("H17").Value = -13.461538461538
$ws.Range("I17").Value = 226
$ws.Range("J17").Value = 239
$ws.Range("K17").Value = -5.439330543933
$ws.Range("L17").Value = 0.444444444444
$ws.Range("M17").Value = 128.282828282828
$ws.Range("N17").Value = -37.222222222222

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -52.173913043478
$ws.Range("I18").Value = 74
$ws.Range("J18").Value = 127
$ws.Range("K18").Value = -41.732283464566
$ws.Range("L18").Value = 4.225352112676
$ws.Range("M18").Value = 51.020408163265
$ws.Range("N18").Value = -81.773399014778

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 23.076923076923
$ws.Range("I19").Value = 168
$ws.Range("J19").Value = 154
$ws.Range("K19").Value = 9.090909090909
$ws.Range("L19").Value = 24.444444444444
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 35.483870967741

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 133.333333333333
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 70
$ws.Range("I20").Value = 179
$ws.Range("J20").Value = 109
$ws.Range("K20").Value = 64.220183486238
$ws.Range("L20").Value = 113.095238095238
$ws.Range("M20").Value = 297.777777777778
$ws.Range("N20").Value = -16.744186046511

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -19.512195121951
$ws.Range("F21").Value = 151
$ws.Range("G21").Value = 148
$ws.Range("H21").Value = 2.027027027027
$ws.Range("I21").Value = 812
$ws.Range("J21").Value = 773
$ws.Range("K21").Value = 5.045278137128
$ws.Range("L21").Value = 25.115562403698
$ws.Range("M21").Value = 110.362694300518
$ws.Range("N21").Value = -48.769716088328

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("M22").Value = -44.444444444444

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -83.333333333333
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 36
$ws.Range("H23").Value = -55.555555555555
$ws.Range("I23").Value = 154
$ws.Range("J23").Value = 137
$ws.Range("K23").Value = 12.408759124087
$ws.Range("L23").Value = 108.108108108108
$ws.Range("M23").Value = 120

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 10.714285714285
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = 3.191489361702
$ws.Range("I24").Value = 449
$ws.Range("J24").Value = 424
$ws.Range("K24").Value = 5.896226415094
$ws.Range("L24").Value = 22.010869565217
$ws.Range("M24").Value = 56.445993031358

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = -37.5
$ws.Range("F25").Value = 85
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = 6.25
$ws.Range("I25").Value = 442
$ws.Range("J25").Value = 389
$ws.Range("K25").Value = 13.624678663239
$ws.Range("L25").Value = 31.157270029673
$ws.Range("M25").Value = 38.125

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
Set-PlaceholderZero "C26"
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = -20
$ws.Range("J26").Value = 26
$ws.Range("K26").Value = 0

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 9
$ws.Range("H27").Value = 80
$ws.Range("I27").Value = 42
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = 61.538461538461
$ws.Range("L27").Value = 100

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------------
$ws.Range("G28").Value = 5
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -28.571428571428
$ws.Range("L28").Value = -56.521739130434
$ws.Range("M28").Value = -50
$ws.Range("N28").Value = -66.666666666666

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Range("G29").Value = 4
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -23.076923076923
$ws.Range("L29").Value = -47.368421052631
$ws.Range("M29").Value = -37.5
$ws.Range("N29").Value = -66.666666666666
